$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("展览")
$ws4 = $wb.Worksheets.Item("全部类型")

# --- Sheet: 展览 ---
$ws1.Cells.Item(2,6).Value = 169
$ws1.Cells.Item(3,6).Value = 174
$ws1.Cells.Item(5,6).Value = 4964
$ws1.Cells.Item(9,6).Value = 543
$ws1.Cells.Item(10,6).Value = 505
$ws1.Cells.Item(11,6).Value = 417
$ws1.Cells.Item(13,6).Value = 1373
$ws1.Cells.Item(14,6).Value = 3602
$ws1.Cells.Item(15,6).Value = 406
$ws1.Cells.Item(16,6).Value = 132
$ws1.Cells.Item(17,6).Value = 115
$ws1.Cells.Item(19,6).Value = 2615
$ws1.Cells.Item(21,3).Value = '南昌·第二届龙年动漫展'
$ws1.Cells.Item(21,4).Value = '南龙潘街666号二楼万达影城斜对面 融创茂啃趣馆'
$ws1.Cells.Item(21,5).Value = '2024.05.03 10:00-05.04 18:00'
$ws1.Cells.Item(21,6).Value = 2
$ws1.Cells.Item(21,7).Value = 55
$ws1.Cells.Item(21,8).Value = 'https://show.bilibili.com/platform/detail.html?id=83954'
$ws1.Cells.Item(21,9).Value = '//i1.hdslb.com/bfs/openplatform/202404/J1ebVrhR1712152650906.jpeg'
$ws1.Cells.Item(22,3).Value = '新余·LD02国风动漫嘉年华'
$ws1.Cells.Item(22,4).Value = '劳动北路888号 金联体育篮球馆'
$ws1.Cells.Item(22,5).Value = '2024.05.03 10:00-05.03 17:00'
$ws1.Cells.Item(22,6).Value = 83
$ws1.Cells.Item(22,7).Value = 50
$ws1.Cells.Item(22,8).Value = 'https://show.bilibili.com/platform/detail.html?id=83379'
$ws1.Cells.Item(22,9).Value = '//i0.hdslb.com/bfs/openplatform/202403/jozduadT1711362183223.jpeg'
$ws1.Cells.Item(23,3).Value = '江西·2024南昌玛雅《次元之芯》主题动漫嘉年华'
$ws1.Cells.Item(23,4).Value = '真君路999号 南昌玛雅乐园'
$ws1.Cells.Item(23,5).Value = '2024.05.03 10:30-05.04 19:30'
$ws1.Cells.Item(23,6).Value = 36
$ws1.Cells.Item(23,7).Value = '不可售'
$ws1.Cells.Item(23,8).Value = 'https://show.bilibili.com/platform/detail.html?id=83608'
$ws1.Cells.Item(23,9).Value = '//i2.hdslb.com/bfs/openplatform/202404/83wvFhen1712040649705.jpeg'
$ws1.Cells.Item(24,3).Value = '江西·ShiningStaR数字互娱嘉年华配音演员史泽鲲专场见面会'
$ws1.Cells.Item(24,4).Value = '前湖大道欣悦湖体育馆 欣悦湖体育馆'
$ws1.Cells.Item(24,5).Value = '2024.05.03 09:30-05.03 17:30'
$ws1.Cells.Item(24,6).Value = 175
$ws1.Cells.Item(24,7).Value = 188
$ws1.Cells.Item(24,8).Value = 'https://show.bilibili.com/platform/detail.html?id=83497'
$ws1.Cells.Item(24,9).Value = '//i1.hdslb.com/bfs/openplatform/202403/qm19B8RF1711620646864.jpeg'
$ws1.Cells.Item(25,2).NumberFormat = "@"
$ws1.Cells.Item(25,2).Value = '2024-05-03'
$ws1.Cells.Item(25,3).Value = '赣州·漫库书店次元漫展'
$ws1.Cells.Item(25,4).Value = '南门口地一大道下沉广场 漫库书店'
$ws1.Cells.Item(25,5).Value = '2024.05.03 10:00-05.04 18:00'
$ws1.Cells.Item(25,6).Value = 43
$ws1.Cells.Item(25,7).Value = 45
$ws1.Cells.Item(25,8).Value = 'https://show.bilibili.com/platform/detail.html?id=83855'
$ws1.Cells.Item(25,9).Value = '//i0.hdslb.com/bfs/openplatform/202404/juDVRy6Y1712481590113.jpeg'
$ws1.Cells.Item(26,2).NumberFormat = "@"
$ws1.Cells.Item(26,2).Value = '2024-05-04'
$ws1.Cells.Item(26,3).Value = '江西·ShiningStaR数字互娱嘉年华 配音演员陈张太康、张惠霖专场见面会'
$ws1.Cells.Item(26,4).Value = '前湖大道欣悦湖体育馆 欣悦湖体育馆'
$ws1.Cells.Item(26,5).Value = '2024.05.04 09:30-05.04 17:30'
$ws1.Cells.Item(26,6).Value = 127
$ws1.Cells.Item(26,7).Value = 228
$ws1.Cells.Item(26,8).Value = 'https://show.bilibili.com/platform/detail.html?id=83593'
$ws1.Cells.Item(26,9).Value = '//i0.hdslb.com/bfs/openplatform/202404/LcnCzDxF1711935576170.jpeg'
$ws1.Cells.Item(27,2).NumberFormat = "@"
$ws1.Cells.Item(27,2).Value = '2024-05-18'
$ws1.Cells.Item(27,3).Value = '南昌·花绒万兽首届兽聚'
$ws1.Cells.Item(27,4).Value = '南龙潘街666号二楼万达影城斜对面 融创茂啃趣馆'
$ws1.Cells.Item(27,5).Value = '2024.05.18 09:30-05.19 16:30'
$ws1.Cells.Item(27,6).Value = 57
$ws1.Cells.Item(27,7).Value = 60
$ws1.Cells.Item(27,8).Value = 'https://show.bilibili.com/platform/detail.html?id=83689'
$ws1.Cells.Item(27,9).Value = '//i2.hdslb.com/bfs/openplatform/202403/h4iL6IvI1711790121140.jpeg'
$ws1.Cells.Item(28,2).NumberFormat = "@"
$ws1.Cells.Item(28,2).Value = '2024-05-26'
$ws1.Cells.Item(28,3).Value = '南昌·代号鸢盛花行only'
$ws1.Cells.Item(28,4).Value = '民德路411号 东方豪景花园酒店(民德路店)'
$ws1.Cells.Item(28,5).Value = '2024.05.26 09:30-05.26 17:30'
$ws1.Cells.Item(28,6).Value = 264
$ws1.Cells.Item(28,7).Value = 78
$ws1.Cells.Item(28,8).Value = 'https://show.bilibili.com/platform/detail.html?id=82529'
$ws1.Cells.Item(28,9).Value = '//i1.hdslb.com/bfs/openplatform/202403/TJ8iC73c1709804909450.png'
$ws1.Cells.Item(29,1).Value = 28
$ws1.Cells.Item(29,2).NumberFormat = "@"
$ws1.Cells.Item(29,2).Value = '2024-06-10'
$ws1.Cells.Item(29,3).Value = '上饶·ETI动漫节'
$ws1.Cells.Item(29,4).Value = '滨江东路与体育馆路交叉口西100米 力加体育综合运动中心'
$ws1.Cells.Item(29,5).Value = '2024.06.10 10:00-06.10 16:00'
$ws1.Cells.Item(29,6).Value = 45
$ws1.Cells.Item(29,7).Value = 36.6
$ws1.Cells.Item(29,8).Value = 'https://show.bilibili.com/platform/detail.html?id=83422'
$ws1.Cells.Item(29,9).Value = '//i1.hdslb.com/bfs/openplatform/202403/vvJKFJal1711460768984.jpeg'

# --- Sheet: 全部类型 ---
$ws4.Cells.Item(2,6).Value = 169
$ws4.Cells.Item(3,6).Value = 174
$ws4.Cells.Item(6,6).Value = 4964
$ws4.Cells.Item(10,6).Value = 543
$ws4.Cells.Item(11,6).Value = 505
$ws4.Cells.Item(12,6).Value = 418
$ws4.Cells.Item(14,6).Value = 1373
$ws4.Cells.Item(15,6).Value = 3602
$ws4.Cells.Item(16,6).Value = 406
$ws4.Cells.Item(17,6).Value = 132
$ws4.Cells.Item(18,6).Value = 115
$ws4.Cells.Item(20,6).Value = 2615
$ws4.Cells.Item(22,3).Value = '南昌·第二届龙年动漫展'
$ws4.Cells.Item(22,4).Value = '南龙潘街666号二楼万达影城斜对面 融创茂啃趣馆'
$ws4.Cells.Item(22,5).Value = '2024.05.03 10:00-05.04 18:00'
$ws4.Cells.Item(22,6).Value = 2
$ws4.Cells.Item(22,7).Value = 55
$ws4.Cells.Item(22,8).Value = 'https://show.bilibili.com/platform/detail.html?id=83954'
$ws4.Cells.Item(22,9).Value = '//i1.hdslb.com/bfs/openplatform/202404/J1ebVrhR1712152650906.jpeg'
$ws4.Cells.Item(23,3).Value = '新余·LD02国风动漫嘉年华'
$ws4.Cells.Item(23,4).Value = '劳动北路888号 金联体育篮球馆'
$ws4.Cells.Item(23,5).Value = '2024.05.03 10:00-05.03 17:00'
$ws4.Cells.Item(23,6).Value = 83
$ws4.Cells.Item(23,7).Value = 50
$ws4.Cells.Item(23,8).Value = 'https://show.bilibili.com/platform/detail.html?id=83379'
$ws4.Cells.Item(23,9).Value = '//i0.hdslb.com/bfs/openplatform/202403/jozduadT1711362183223.jpeg'
$ws4.Cells.Item(24,3).Value = '江西·2024南昌玛雅《次元之芯》主题动漫嘉年华'
$ws4.Cells.Item(24,4).Value = '真君路999号 南昌玛雅乐园'
$ws4.Cells.Item(24,5).Value = '2024.05.03 10:30-05.04 19:30'
$ws4.Cells.Item(24,6).Value = 36
$ws4.Cells.Item(24,7).Value = '不可售'
$ws4.Cells.Item(24,8).Value = 'https://show.bilibili.com/platform/detail.html?id=83608'
$ws4.Cells.Item(24,9).Value = '//i2.hdslb.com/bfs/openplatform/202404/83wvFhen1712040649705.jpeg'
$ws4.Cells.Item(25,3).Value = '江西·ShiningStaR数字互娱嘉年华配音演员史泽鲲专场见面会'
$ws4.Cells.Item(25,4).Value = '前湖大道欣悦湖体育馆 欣悦湖体育馆'
$ws4.Cells.Item(25,5).Value = '2024.05.03 09:30-05.03 17:30'
$ws4.Cells.Item(25,6).Value = 175
$ws4.Cells.Item(25,7).Value = 188
$ws4.Cells.Item(25,8).Value = 'https://show.bilibili.com/platform/detail.html?id=83497'
$ws4.Cells.Item(25,9).Value = '//i1.hdslb.com/bfs/openplatform/202403/qm19B8RF1711620646864.jpeg'
$ws4.Cells.Item(26,2).NumberFormat = "@"
$ws4.Cells.Item(26,2).Value = '2024-05-03'
$ws4.Cells.Item(26,3).Value = '赣州·漫库书店次元漫展'
$ws4.Cells.Item(26,4).Value = '南门口地一大道下沉广场 漫库书店'
$ws4.Cells.Item(26,5).Value = '2024.05.03 10:00-05.04 18:00'
$ws4.Cells.Item(26,6).Value = 43
$ws4.Cells.Item(26,7).Value = 45
$ws4.Cells.Item(26,8).Value = 'https://show.bilibili.com/platform/detail.html?id=83855'
$ws4.Cells.Item(26,9).Value = '//i0.hdslb.com/bfs/openplatform/202404/juDVRy6Y1712481590113.jpeg'
$ws4.Cells.Item(27,2).NumberFormat = "@"
$ws4.Cells.Item(27,2).Value = '2024-05-04'
$ws4.Cells.Item(27,3).Value = '江西·ShiningStaR数字互娱嘉年华 配音演员陈张太康、张惠霖专场见面会'
$ws4.Cells.Item(27,4).Value = '前湖大道欣悦湖体育馆 欣悦湖体育馆'
$ws4.Cells.Item(27,5).Value = '2024.05.04 09:30-05.04 17:30'
$ws4.Cells.Item(27,6).Value = 127
$ws4.Cells.Item(27,7).Value = 228
$ws4.Cells.Item(27,8).Value = 'https://show.bilibili.com/platform/detail.html?id=83593'
$ws4.Cells.Item(27,9).Value = '//i0.hdslb.com/bfs/openplatform/202404/LcnCzDxF1711935576170.jpeg'
$ws4.Cells.Item(28,2).NumberFormat = "@"
$ws4.Cells.Item(28,2).Value = '2024-05-18'
$ws4.Cells.Item(28,3).Value = '南昌·花绒万兽首届兽聚'
$ws4.Cells.Item(28,4).Value = '南龙潘街666号二楼万达影城斜对面 融创茂啃趣馆'
$ws4.Cells.Item(28,5).Value = '2024.05.18 09:30-05.19 16:30'
$ws4.Cells.Item(28,6).Value = 57
$ws4.Cells.Item(28,7).Value = 60
$ws4.Cells.Item(28,8).Value = 'https://show.bilibili.com/platform/detail.html?id=83689'
$ws4.Cells.Item(28,9).Value = '//i2.hdslb.com/bfs/openplatform/202403/h4iL6IvI1711790121140.jpeg'
$ws4.Cells.Item(29,2).NumberFormat = "@"
$ws4.Cells.Item(29,2).Value = '2024-05-26'
$ws4.Cells.Item(29,3).Value = '南昌·代号鸢盛花行only'
$ws4.Cells.Item(29,4).Value = '民德路411号 东方豪景花园酒店(民德路店)'
$ws4.Cells.Item(29,5).Value = '2024.05.26 09:30-05.26 17:30'
$ws4.Cells.Item(29,6).Value = 264
$ws4.Cells.Item(29,7).Value = 78
$ws4.Cells.Item(29,8).Value = 'https://show.bilibili.com/platform/detail.html?id=82529'
$ws4.Cells.Item(29,9).Value = '//i1.hdslb.com/bfs/openplatform/202403/TJ8iC73c1709804909450.png'
$ws4.Cells.Item(30,1).Value = 29
$ws4.Cells.Item(30,2).NumberFormat = "@"
$ws4.Cells.Item(30,2).Value = '2024-06-10'
$ws4.Cells.Item(30,3).Value = '上饶·ETI动漫节'
$ws4.Cells.Item(30,4).Value = '滨江东路与体育馆路交叉口西100米 力加体育综合运动中心'
$ws4.Cells.Item(30,5).Value = '2024.06.10 10:00-06.10 16:00'
$ws4.Cells.Item(30,6).Value = 45
$ws4.Cells.Item(30,7).Value = 36.6
$ws4.Cells.Item(30,8).Value = 'https://show.bilibili.com/platform/detail.html?id=83422'
$ws4.Cells.Item(30,9).Value = '//i1.hdslb.com/bfs/openplatform/202403/vvJKFJal1711460768984.jpeg'
